# Apply updated cryptocurrency price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '54.593.36'
$ws.Cells.Item(2, 5).Value = '  +0.54%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.281.68'
$ws.Cells.Item(3, 5).Value = '  +0.01%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.10%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''505.67'
$ws.Cells.Item(5, 5).Value = '  +1.38%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''128.79'
$ws.Cells.Item(6, 5).Value = '  +0.09%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.22%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '''0.528'
$ws.Cells.Item(8, 5).Value = '  +0.01%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '2.301.07'
$ws.Cells.Item(9, 5).Value = '  +0.57%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''0.0969'
$ws.Cells.Item(10, 5).Value = '  +1.56%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +1.57%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +2.24%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +4.12%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '''23.44'
$ws.Cells.Item(14, 5).Value = '  +3.04%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '2.689.71'
$ws.Cells.Item(15, 5).Value = '  +0.18%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '54.637.48'
$ws.Cells.Item(16, 5).Value = '  +0.77%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  +1.41%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '2.294.79'
$ws.Cells.Item(18, 5).Value = '  -0.30%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +1.53%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  +1.05%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''307.60'
$ws.Cells.Item(21, 5).Value = '  +1.01%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''6.59'
$ws.Cells.Item(22, 5).Value = '  +2.79%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''0.999'
$ws.Cells.Item(23, 5).Value = '  -0.08%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''60.41'
$ws.Cells.Item(24, 5).Value = '  -2.31%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''0.994'
$ws.Cells.Item(25, 5).Value = '  -0.61%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +2.21%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '''171.53'
$ws.Cells.Item(28, 5).Value = '  -1.87%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +2.06%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +2.56%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +0.89%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +4.15%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  -0.01%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '''17.99'
$ws.Cells.Item(34, 5).Value = '  +1.15%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''0.995'
$ws.Cells.Item(35, 5).Value = '  -0.03%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''0.909'
$ws.Cells.Item(36, 5).Value = '  -2.66%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''1.21'
$ws.Cells.Item(37, 5).Value = '  +0.81%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +1.81%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '''36.54'
$ws.Cells.Item(39, 5).Value = '  +1.33%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +0.42%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +1.02%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'RenderToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(42, 4).Value = '''5.06'
$ws.Cells.Item(42, 5).Value = '  +5.94%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''131.43'
$ws.Cells.Item(43, 5).Value = '  +5.21%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'Filecoin'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(44, 4).Value = '''3.41'
$ws.Cells.Item(44, 5).Value = '  +0.79%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '''251.41'
$ws.Cells.Item(45, 5).Value = '  +4.73%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +1.02%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '''0.0910'
$ws.Cells.Item(47, 5).Value = '  +1.64%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''0.551'
$ws.Cells.Item(48, 5).Value = '  +0.98%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +0.37%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +0.63%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +0.37%  '
